# "create material, salary as manager"
# Append three new data rows (10-12) to Sheet1, extending the table that
# currently ends at row 9 (dimension A1:E9 -> A1:E12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A on existing data rows carries style index 1 (bold, thin border,
# centered/top-aligned, numeric). Copy that formatting onto the new rows'
# A cells before writing their values so we reuse the existing style
# rather than minting a new one.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)  # xlPasteFormats

# Row 10: 8 | dwdwdw | суммы | 12 | effrfwwpkp
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "dwdwdw"
$ws.Range("C10").Value = "суммы"
$ws.Range("E10").Value = "effrfwwpkp"

# Row 11: 9 | frwrfr | суммы | 8 | effrfwwpkp
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "frwrfr"
$ws.Range("C11").Value = "суммы"
$ws.Range("E11").Value = "effrfwwpkp"

# Row 12: 10 | heyy | доллары | 12 | effrfwwpkp
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "heyy"
$ws.Range("C12").Value = "доллары"
$ws.Range("E12").Value = "effrfwwpkp"

# Column D holds numeric-looking values but must stay text (like the rest
# of column D in the existing rows), so force text formatting, write the
# values, then drop back to the unstyled "Normal" cell style (matching
# the source, which has no explicit style on these cells).
$ws.Range("D10:D12").NumberFormat = "@"
$ws.Range("D10").Value = "12"
$ws.Range("D11").Value = "8"
$ws.Range("D12").Value = "12"
$ws.Range("D10:D12").Style = "Normal"
